# Update MSME Country Indicators - Timor-Leste Summary values.
# Row 12: "Enterprises (absolute #)"          -> Micro / SMEs / MSMEs
# Row 13: "Enterprises density (per 1000 people)" -> Micro / SMEs / MSMEs
#
# These cells store text (shared-string) values such as "3008.2", not
# numbers, so the replacement values must also be written as text. Excel
# normally auto-converts a numeric-looking string into a real number when
# assigned via .Value, which would both change the stored type and drop
# the cell's existing style. To avoid that, the cell's style is captured,
# the format is switched to Text ("@") just long enough to assign the
# value as a string, and then the original style is restored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$cells  = @("B12", "C12", "D12", "B13", "C13", "D13")
$values = @("3008.22", "1130.16", "4138.38", "3.11", "1.17", "4.28")

for ($i = 0; $i -lt $cells.Length; $i++) {
    $cell = $ws.Range($cells[$i])
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $values[$i]
    $cell.Style = $originalStyle
}
